# Updates the Price (D) and Volume(1h) (E) columns of the crypto symbol
# list with refreshed quotes, matching the GitHub Actions data-refresh
# commit. Values are text-formatted (not numeric) in the source sheet,
# so each assignment uses a leading apostrophe to force Excel to keep
# them as text (preserving exact formatting such as trailing zeros),
# then resets the cell style to Normal so no stray number-format is left
# behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'303.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'5.09%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'34.80"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.97%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.168"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'4.31%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'6.35%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-2.43%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'8.049"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'4.18%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'3.993"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'7.18%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9216"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'1.11%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'8.30%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1840"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'7.88%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08487"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'3.81%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03388"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'8.81%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09918"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.61%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001483"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.20%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.04654"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'2.96%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005791"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.97%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'3.470"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.01%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'1.65%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'2.76%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1328"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.14%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'4.549"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.27%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.2400"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'14.20%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001223"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.91%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004327"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'3.68%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.04%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'0.11%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01734"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'9.59%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04752"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'6.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007690"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.41%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1412"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'6.10%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.007261"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-26.56%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002219"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.01%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.009963"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'13.39%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006068"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.64%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'51.24%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002695"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'34.68%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("E51").Style = "Normal"
